$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56, shifting rows 56-58 down to 57-59
$ws.Rows.Item(56).Insert()

# Fill in the new row 56 with the new data point
$ws.Cells.Item(56, 1).Value = 10
$ws.Cells.Item(56, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(56, 3).Value = "La Araucanía"
$ws.Cells.Item(56, 4).Value = 44516
$ws.Cells.Item(56, 4).NumberFormat = $ws.Cells.Item(57, 4).NumberFormat
$ws.Cells.Item(56, 5).Value = 9
$ws.Cells.Item(56, 6).Value = "Fruta"
$ws.Cells.Item(56, 7).Value = 100101
$ws.Cells.Item(56, 8).Value = "Berries"
$ws.Cells.Item(56, 9).Value = 100101001
$ws.Cells.Item(56, 10).Value = "Arándano (blue)"
$ws.Cells.Item(56, 11).Value = "Sin especificar"
$ws.Cells.Item(56, 12).Value = "Primera"
$ws.Cells.Item(56, 13).Value = 530
$ws.Cells.Item(56, 14).Value = 3300
$ws.Cells.Item(56, 15).Value = 3500
$ws.Cells.Item(56, 16).Value = 3432
$ws.Cells.Item(56, 17).Value = "$/kilo"
$ws.Cells.Item(56, 18).Value = "Región del Maule"
$ws.Cells.Item(56, 19).Value = 3432
$ws.Cells.Item(56, 20).Value = 1

# Now fix row 57 (previously row 56, shifted down by the insert): update R value
$ws.Cells.Item(57, 18).Value = "Región del Maule"
